# Update countries & provincias Spain
# Applies the 11-Jul-2020 12:37 data refresh to the "Pais" sheet:
#   - Updates the "last updated" timestamp string
#   - Re-sorts a handful of countries whose case counts changed rank
#     (country names move, numeric columns B-H get refreshed)
#   - Refreshes case totals (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) for the
#     countries whose figures changed, without touching rows that were
#     not updated in the source feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# ---------------------------------------------------------------
# 1) Update the "Datos actualizados..." timestamp banner (A1)
# ---------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 11 de Julio de 2020 a las 12:37"

# ---------------------------------------------------------------
# 2) Countries that swapped rank (country name moves to the other
#    row, each row keeps the new stats that belong to the country
#    now occupying that row).
# ---------------------------------------------------------------
$ws.Range("A37").Value = "Kuwait"
$ws.Range("A38").Value = "Emiratos Arabes Unidos"

$ws.Range("A50").Value = "Rumania"
$ws.Range("A51").Value = "Barein"
$ws.Range("A52").Value = "Armenia"

$ws.Range("A80").Value = "Senegal"
$ws.Range("A81").Value = "Consejo Danes para los Refugiados"

$ws.Range("A92").Value = "Estado de Palestina"
$ws.Range("A93").Value = "Guayana Francesa"

# ---------------------------------------------------------------
# 3) Refresh numeric columns (B=Casos totales, C=Nuevos casos,
#    D=Casos activos, E=Recuperados, F=Casos criticos,
#    G=Muertes hoy, H=Muertes) for every updated row.
# ---------------------------------------------------------------
$rowData = @{
    4   = @(3292257,  471, 1460644, 1694931, 0, 11, 136682)
    6   = @( 823927, 1324,  516338,  285436, 0,  9,  22153)
    13  = @( 255117, 2397,  217666,   24816, 0, 188, 12635)
    37  = @(  54058,  478,   43961,    9711, 0,  3,    386)
    38  = @(  54050,    0,   43969,    9751, 0,  0,    330)
    49  = @(  32798,  108,   29400,    1432, 0,  0,   1966)
    50  = @(  32079,  698,   21414,    8794, 0, 24,   1871)
    51  = @(  32039,    0,   27213,    4722, 0,  0,    104)
    52  = @(  31392,  489,   19419,   11414, 0, 13,    559)
    65  = @(  15464,  136,   11895,    3325, 0,  1,    244)
    76  = @(   9391,  249,    5516,    3621, 0,  5,    254)
    79  = @(   8704,    8,    8515,      67, 0,  1,    122)
    80  = @(   8014,  132,    5381,    2488, 0,  0,    145)
    81  = @(   7905,    0,    3513,    4203, 0,  0,    189)
    92  = @(   5931,  380,     536,    5366, 0,  2,     29)
    93  = @(   5704,    0,    2785,    2896, 0,  0,     23)
    102 = @(   3371,   93,    1881,    1401, 0,  4,     89)
    158 = @(    462,    4,     118,     321, 0,  0,     23)
}

$cols = @("B", "C", "D", "E", "F", "G", "H")

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}
